# Weekly update: insert a new record (week of 2022-06-24) above the
# existing history for "Bruselas (repollito)" at Mercado Mayorista Lo
# Valledor de Santiago, pushing prior rows 16-19 down to 17-20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 16; existing rows 16-19 shift to 17-20.
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the new weekly data point.
$ws.Range("A16").Value = 6
$ws.Range("B16").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C16").Value = "Metropolitana"
$ws.Range("D16").Value = 44736
$ws.Range("E16").Value = 13
$ws.Range("F16").Value = 100112035
$ws.Range("G16").Value = "Bruselas (repollito)"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 180
$ws.Range("K16").Value = 17000
$ws.Range("L16").Value = 19000
$ws.Range("M16").Value = 17889
$ws.Range("N16").Value = "$/malla 15 kilos"
$ws.Range("O16").Value = "Provincia de Quillota"
$ws.Range("P16").Value = 1193
$ws.Range("Q16").Value = 15
$ws.Range("R16").Value = "Hortaliza"
